$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column A on the new rows to be treated as plain text (same as the
# existing date-like strings, e.g. A27 "03-08-2021"), otherwise Excel will
# auto-recognize the "dd-mm-yyyy" pattern and convert it into a date serial
# number.
$ws.Range("A28:A29").NumberFormat = "@"

# New row 28: 04-08-2021
$ws.Range("A28").Value = "04-08-2021"
$ws.Range("B28").Value = 200000
$ws.Range("C28").Value = 376000
$ws.Range("D28").Value = 100000
$ws.Range("E28").Value = 74000
$ws.Range("F28").Value = 26000
$ws.Range("G28").Value = 0.85

# New row 29: 05-08-2021
$ws.Range("A29").Value = "05-08-2021"
$ws.Range("B29").Value = 200000
$ws.Range("C29").Value = 281000
$ws.Range("D29").Value = 100000
$ws.Range("E29").Value = 40000
$ws.Range("F29").Value = 60000
$ws.Range("G29").Value = 0.95

# Reapplying the default "Normal" style clears the explicit style index
# that NumberFormat left behind, so A28:A29 end up unstyled (default),
# matching every other data cell in column A (e.g. A27).
$ws.Range("A28:A29").Style = "Normal"
